$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are plain text (e.g. "27.743.99") in the source data.
# Forcing NumberFormat to Text ("@") before assignment stops Excel from
# re-interpreting digit-and-dot strings as numbers; resetting the Style
# back to "Normal" afterwards restores the original (unstyled) cell format
# so only the cell VALUE changes, matching the source diff.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.743.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.55%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.639.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.52%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("E6").Value = "  -1.87%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.25"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.46%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.262"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.14%  "

$ws.Range("E10").Value = "  +0.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0888"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.02%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.870.84"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.637.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.60%  "

$ws.Range("E14").Value = "  +0.49%  "

$ws.Range("E15").Value = "  -3.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.691.12"
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.93%  "

$ws.Range("E20").Value = "  -0.05%  "

$ws.Range("E22").Value = "  -0.32%  "

$ws.Range("E23").Value = "  +4.59%  "

$ws.Range("E24").Value = "  +3.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.68%  "

$ws.Range("E27").Value = "  -0.84%  "

$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("E31").Value = "  +0.22%  "

$ws.Range("E32").Value = "  +0.13%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.466.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.76%  "

$ws.Range("E35").Value = "  -2.61%  "

$ws.Range("E36").Value = "  -0.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.570"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("E38").Value = "  -0.09%  "

$ws.Range("E39").Value = "  +0.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.902"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.76%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "68.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.77%  "

$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("E43").Value = "  -1.75%  "

$ws.Range("E44").Value = "  +1.45%  "

$ws.Range("E45").Value = "  -1.16%  "

$ws.Range("E46").Value = "  -0.69%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.780.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0106"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.97%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0994"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.11%  "
